$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 1651
$ws1.Range("F11").Value = 1530
$ws1.Range("F13").Value = 48
$ws1.Range("F14").Value = 377
$ws1.Range("F16").Value = 192
$ws1.Range("F18").Value = 21
$ws1.Range("F21").Value = 273

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 1651
$ws4.Range("F12").Value = 1530
$ws4.Range("F14").Value = 48
$ws4.Range("F15").Value = 377
$ws4.Range("F17").Value = 192
$ws4.Range("F19").Value = 21
$ws4.Range("F22").Value = 273
